$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.386024832725525
$ws.Range("B1").Value = 1.462569952011108
$ws.Range("C1").Value = 1.641091227531433
$ws.Range("D1").Value = 2.522682666778564
$ws.Range("E1").Value = 15
